$d = $word.ActiveDocument

# 1. Merge the split runs in the first paragraph ("(Línea modificada) " + bookmark + "!Saludos!")
#    into a single run "(Línea modificada) !Saludos!", dropping the _GoBack bookmark that used
#    to sit between them (it gets re-created at the new rollback note below).
$d.Content.Find.Execute("modificada) !Saludos!", $true, $false, $false, $false, $false, $true, 1, $false, "modificada) !Saludos!", 2)

# 2. Append a new blank paragraph, then a paragraph with the rollback note.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()

$blank = $d.Paragraphs($d.Paragraphs.Count)
$blank.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$insertRange = $newPara.Range
# Trailing "ZZ" placeholder keeps the bookmark insertion point away from the paragraph-mark
# edge (an edge Range resolves incorrectly); it is stripped again right after.
$insertRange.InsertBefore("Mireia ha subido su carpeta. Volvemos al commit de MarvinZZ")

# 3. Put the _GoBack bookmark right after the new sentence, matching where it used to live.
$bookmarkSpot = $d.Content
$bookmarkSpot.Find.Execute("commit de Marvin") | Out-Null
$bookmarkSpot.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

# 4. Remove the temporary placeholder.
$placeholder = $d.Content
$placeholder.Find.Execute("ZZ") | Out-Null
$placeholder.Delete()
